# CHG: Docu week1 - shortened "zum Beispiel" to "z.B."
#
# This reproduces, via the Word COM/OM surface, the content edit shown
# in the unified diff:
#   - " (zum Beispiel " -> " (" + "z.B." + " " (three separate runs,
#     matching the run layout produced by Word when text with different
#     history/formatting boundaries is edited)
#   - the stray "_GoBack" bookmark that used to sit at the end of the
#     "Graphikprogramm ... Spielfelds" bullet is gone, and a fresh
#     "_GoBack" bookmark now sits right after the new "z.B. " text
#     (immediately before "Smartphone") - exactly where Word leaves the
#     go-back mark after the last interactive edit.

$d = $word.ActiveDocument

# --- 1. Drop the old "_GoBack" bookmark (end of the Graphikprogramm bullet) ---
try {
    $oldGoBack = $d.Bookmarks.Item("_GoBack")
    $oldGoBack.Delete()
} catch {
    # no pre-existing _GoBack bookmark - nothing to remove
}

# --- 2. Locate the sentence fragment we need to rework ---
# "...Geräten (zum Beispiel Smartphones) spielbar sein."
$found = $d.Content
$found.Find.Execute("en (zum Beispiel Smartphone", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base = $found.Start

# Character offsets (relative to $base) inside "en (zum Beispiel Smartphone":
#   0123456789...
#   e n   ( z  u  m     B  e  i  s  p  i  e  l     S  m  a ...
#   0 1 2 3 4  5  6  7  8  9  10 11 12 13 14 15 16 17 ...
$posEnBoundary    = $base + 2    # between "en" and " ("
$posParenBoundary = $base + 4    # between "(" and "zum Beispiel "
$zumStart         = $base + 4    # start of "zum Beispiel"
$zumEnd           = $base + 16   # end of "zum Beispiel" (before the trailing space)

# --- 3. Pin down the future run boundaries with temporary bookmarks ---
# Adding a bookmark at a position forces a clean run split there (with no
# left-over run-formatting), which is what lets us reproduce the exact
# three-run layout from the diff instead of Word merging everything back
# into a single run.
$d.Bookmarks.Add("ZZ_TMP_EN", $d.Range($posEnBoundary, $posEnBoundary)) | Out-Null
$d.Bookmarks.Add("ZZ_TMP_PAREN", $d.Range($posParenBoundary, $posParenBoundary)) | Out-Null
$d.Bookmarks.Add("ZZ_TMP_SPACE", $d.Range($zumEnd, $zumEnd)) | Out-Null

# --- 4. Now that the run is fully isolated on both sides, shorten the text ---
$target = $d.Range($zumStart, $zumEnd)
$target.Text = "z.B."

# --- 5. Re-find the new text to drop the fresh "_GoBack" bookmark right
#        before "Smartphone" (after the "z.B. " run) ---
$found2 = $d.Content
$found2.Find.Execute("z.B. Smartphone", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$posGoBack = $found2.Start + 5   # right after "z.B. "
$d.Bookmarks.Add("_GoBack", $d.Range($posGoBack, $posGoBack)) | Out-Null

# --- 6. Clean up the temporary helper bookmarks (the run splits persist) ---
$d.Bookmarks.Item("ZZ_TMP_EN").Delete()
$d.Bookmarks.Item("ZZ_TMP_PAREN").Delete()
$d.Bookmarks.Item("ZZ_TMP_SPACE").Delete()
